$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.576.92"
$ws.Range("E2").Value = "  -7.36%  "
$ws.Range("D3").Value = "1.698.53"
$ws.Range("E3").Value = "  -5.85%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'220.07"
$ws.Range("E5").Value = "  -5.25%  "
$ws.Range("D6").Value = "'0.5130"
$ws.Range("E6").Value = "  -13.40%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -4.78%  "
$ws.Range("D10").Value = "'0.06280"
$ws.Range("E10").Value = "  -7.87%  "
$ws.Range("D11").Value = "'0.07355"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("D12").Value = "1.700.36"
$ws.Range("E12").Value = "  -5.76%  "
$ws.Range("D13").Value = "'4.530"
$ws.Range("E13").Value = "  -4.84%  "
$ws.Range("D14").Value = "'0.5848"
$ws.Range("E14").Value = "  -5.96%  "
$ws.Range("D15").Value = "1.929.52"
$ws.Range("E15").Value = "  -5.83%  "
$ws.Range("D16").Value = "'0.000008451"
$ws.Range("E16").Value = "  -8.28%  "
$ws.Range("D17").Value = "'65.69"
$ws.Range("E17").Value = "  -13.08%  "
$ws.Range("D18").Value = "26.639.71"
$ws.Range("E18").Value = "  -7.02%  "
$ws.Range("D19").Value = "'5.028"
$ws.Range("E19").Value = "  -8.31%  "
$ws.Range("D21").Value = "'10.98"
$ws.Range("E21").Value = "  -4.56%  "
$ws.Range("D22").Value = "'187.30"
$ws.Range("E22").Value = "  -11.08%  "
$ws.Range("E23").Value = "  -8.05%  "
$ws.Range("D25").Value = "'144.90"
$ws.Range("E25").Value = "  -5.79%  "
$ws.Range("D26").Value = "'7.568"
$ws.Range("E26").Value = "  -3.70%  "
$ws.Range("D27").Value = "'0.1154"
$ws.Range("E27").Value = "  -8.92%  "
$ws.Range("D28").Value = "'15.71"
$ws.Range("E28").Value = "  -4.35%  "
$ws.Range("D29").Value = "'1.325"
$ws.Range("E29").Value = "  -7.11%  "
$ws.Range("D30").Value = "'0.05713"
$ws.Range("E30").Value = "  -7.58%  "
$ws.Range("D31").Value = "'1.335"
$ws.Range("E31").Value = "  -6.41%  "
$ws.Range("D32").Value = "'3.530"
$ws.Range("E32").Value = "  -6.70%  "
$ws.Range("D33").Value = "'3.510"
$ws.Range("E33").Value = "  -6.20%  "
$ws.Range("D34").Value = "'1.649"
$ws.Range("E34").Value = "  -5.03%  "
$ws.Range("D35").Value = "'1.030"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("D36").Value = "'0.6033"
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("D37").Value = "'2.373"
$ws.Range("E37").Value = "  -5.10%  "
$ws.Range("D38").Value = "'2.684"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").Value = "1.102.68"
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").Value = "'0.01610"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").Value = "'0.8621"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("D42").Value = "'5.848"
$ws.Range("E42").Value = "  -10.57%  "
$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "'99.04"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").Value = "1.856.36"
$ws.Range("E45").Value = "  -5.18%  "
$ws.Range("D46").Value = "'0.00000000110"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("D47").Value = "'56.70"
$ws.Range("E47").Value = "  -6.41%  "
$ws.Range("D48").Value = "'8.169"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'0.05248"
$ws.Range("E50").Value = "  -4.16%  "
$ws.Range("D51").Value = "'0.4328"
$ws.Range("E51").Value = "  -3.36%  "
